# Adds a "WithoutEndDate" column (header + True/False values) to the
# "AddressTypes Data" sheet, per the commit:
# "Added ability for AddressTypesType withoutEndDate"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddressTypes Data")

$ws.Cells.Item(1, 3).Value = "WithoutEndDate"

# Use a leading apostrophe so Excel stores these as literal text ("True"/
# "False") in the shared-string table instead of native booleans, then
# strip the resulting quote-prefix formatting so the cell style stays the
# default (s="0"), matching the rest of the sheet.
$ws.Cells.Item(2, 3).Value = "'True"
$ws.Cells.Item(2, 3).ClearFormats()

$ws.Cells.Item(3, 3).Value = "'False"
$ws.Cells.Item(3, 3).ClearFormats()

$ws.Cells.Item(4, 3).Value = "'True"
$ws.Cells.Item(4, 3).ClearFormats()
